$d = $word.ActiveDocument
$d.Content.Find.Execute("Fecha: dd/mm/aaaa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fecha: 15/06/2021", 2)
